$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = 47
$ws.Range("G6").Value = 1404.36
$ws.Range("B10").Value = 26195.16
$ws.Range("F59").Value = 17
$ws.Range("G59").Value = 1396.38
$ws.Range("F62").Value = 33
$ws.Range("G62").Value = 2340.36
$ws.Range("F64").Value = 109
$ws.Range("G64").Value = 8847.530000000001
$ws.Range("F86").Value = 44
$ws.Range("G86").Value = 5520.68
$ws.Range("B90").Value = 164339.2
$ws.Range("B127").Value = 57552
$ws.Range("E127").Value = 136.86
$ws.Range("F127").Value = -5
$ws.Range("G127").Value = -603.45
$ws.Range("B128").Value = 64329
$ws.Range("E128").Value = 128.32
$ws.Range("F128").Value = 1
$ws.Range("G128").Value = 120.69
$ws.Range("F151").Value = 86
$ws.Range("G151").Value = 7471.68
$ws.Range("B156").Value = 28321.35
$ws.Range("F192").Value = 0
$ws.Range("G192").Value = 0
$ws.Range("F199").Value = 18
$ws.Range("G199").Value = 4452.84
$ws.Range("F204").Value = 0
$ws.Range("G204").Value = 0
$ws.Range("B216").Value = 30644.53
$ws.Range("F222").Value = 8
$ws.Range("G222").Value = 1159.44
$ws.Range("F249").Value = 132
$ws.Range("G249").Value = 18192.24
$ws.Range("B260").Value = 164575.14
$ws.Range("F294").Value = 21
$ws.Range("G294").Value = 1498.56
$ws.Range("B304").Value = 160679.64
$ws.Range("F338").Value = 70
$ws.Range("G338").Value = 1659
$ws.Range("B346").Value = 23239.08
$ws.Range("F353").Value = 10
$ws.Range("G353").Value = 1371.9
$ws.Range("B358").Value = 34378.99
$ws.Range("B364").Value = 53602
$ws.Range("E364").Value = 15.69
$ws.Range("F364").Value = -231
$ws.Range("G364").Value = -3037.65
$ws.Range("B365").Value = 65068
$ws.Range("E365").Value = 13.97
$ws.Range("F365").Value = 63
$ws.Range("G365").Value = 828.45
$ws.Range("B366").Value = 65066
$ws.Range("E366").Value = 13.61
$ws.Range("F366").Value = 90
$ws.Range("G366").Value = 1152.9
$ws.Range("B367").Value = 53263
$ws.Range("E367").Value = 15.29
$ws.Range("F367").Value = -309
$ws.Range("G367").Value = -3958.29
$ws.Range("B372").Value = 45706
$ws.Range("E372").Value = 23.58
$ws.Range("F372").Value = -202
$ws.Range("G372").Value = -3985.46
$ws.Range("B373").Value = 64922
$ws.Range("E373").Value = 20.98
$ws.Range("F373").Value = 67
$ws.Range("G373").Value = 1321.91
$ws.Range("B380").Value = 64925
$ws.Range("E380").Value = 13.97
$ws.Range("F380").Value = 111
$ws.Range("G380").Value = 1459.65
$ws.Range("B381").Value = 45709
$ws.Range("E381").Value = 15.69
$ws.Range("F381").Value = -300
$ws.Range("G381").Value = -3945
$ws.Range("B385").Value = 53595
$ws.Range("E385").Value = 17.61
$ws.Range("F385").Value = -335
$ws.Range("G385").Value = -4934.55
$ws.Range("B386").Value = 65067
$ws.Range("E386").Value = 15.65
$ws.Range("F386").Value = 126
$ws.Range("G386").Value = 1855.98
$ws.Range("B473").Value = 64830
$ws.Range("E473").Value = 34.9
$ws.Range("F473").Value = 104
$ws.Range("G473").Value = 3414.32
$ws.Range("B474").Value = 60022
$ws.Range("E474").Value = 37.22
$ws.Range("F474").Value = -113
$ws.Range("G474").Value = -3709.79
$ws.Range("B572").Value = 65079
$ws.Range("F572").Value = 6
$ws.Range("G572").Value = 245.22
$ws.Range("B573").Value = 65362
$ws.Range("F573").Value = 18
$ws.Range("G573").Value = 735.66
$ws.Range("F582").Value = 20
$ws.Range("G582").Value = 1139.8
$ws.Range("B583").Value = 12197.78
$ws.Range("F599").Value = 1253
$ws.Range("G599").Value = 204376.83
$ws.Range("B606").Value = 349196.24
$ws.Range("B619").Value = 1540337.65
$ws.Range("B620").Value = 1540337.65
